$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 4146.1816
$ws.Range("I69").Value = 4702.6
$ws.Range("J69").Value = 3682.5
$ws.Range("K69").Value = 14107.8
$ws.Range("L69").Value = 11047.5
$ws.Range("M69").Value = -13233.8
$ws.Range("N69").Value = -12795.5

# Row 72
$ws.Range("H72").Value = 4146.1816
$ws.Range("I72").Value = 4702.6
$ws.Range("J72").Value = 3682.5
$ws.Range("K72").Value = 42323.4
$ws.Range("L72").Value = 33142.5
$ws.Range("M72").Value = -37955.4
$ws.Range("N72").Value = -41878.5

# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

# Row 137
$ws.Range("H137").Value = 14287759
$ws.Range("I137").Value = 1466.6666
$ws.Range("J137").Value = 18184020
$ws.Range("K137").Value = 4399.9998
$ws.Range("L137").Value = 54552060
$ws.Range("M137").Value = -1849.9998
$ws.Range("N137").Value = -54557160

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 19234312
$ws.Range("I74").Value = 29413426
$ws.Range("K74").Value = 29413426
$ws.Range("M74").Value = -29412552

# Row 77
$ws.Range("H77").Value = 19234312
$ws.Range("I77").Value = 29413426
$ws.Range("K77").Value = 147067130
$ws.Range("M77").Value = -147062762

# Row 122
$ws.Range("H122").Value = 8139
$ws.Range("I122").Value = 8944.571
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 26833.713
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -24383.713
$ws.Range("N122").Value = -12400

# Row 138
$ws.Range("H138").Value = 59559.383
$ws.Range("J138").Value = 59559.383
$ws.Range("L138").Value = 59559.383
$ws.Range("N138").Value = -69839.383

$ws = $wb.Worksheets.Item("BSM")
# Row 38
$ws.Range("H38").Value = 20036
$ws.Range("J38").Value = 20036
$ws.Range("L38").Value = 20036
$ws.Range("N38").Value = -20868

# Row 51
$ws.Range("H51").Value = 50780
$ws.Range("J51").Value = 50780
$ws.Range("L51").Value = 50780
$ws.Range("N51").Value = -51762

# Row 86
$ws.Range("H86").Value = 35716264
$ws.Range("I86").Value = 1971.4286
$ws.Range("J86").Value = 71430560
$ws.Range("K86").Value = 1971.4286
$ws.Range("L86").Value = 71430560
$ws.Range("M86").Value = -848.4286
$ws.Range("N86").Value = -71432806

# Row 89
$ws.Range("H89").Value = 35716264
$ws.Range("I89").Value = 1971.4286
$ws.Range("J89").Value = 71430560
$ws.Range("K89").Value = 9857.143
$ws.Range("L89").Value = 357152800
$ws.Range("M89").Value = -4241.143
$ws.Range("N89").Value = -357164032

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 17550816
$ws.Range("I31").Value = 7341.278
$ws.Range("J31").Value = 333333340
$ws.Range("K31").Value = 7341.278
$ws.Range("L31").Value = 333333340
$ws.Range("M31").Value = -7046.278
$ws.Range("N31").Value = -333333930

# Row 34
$ws.Range("H34").Value = 17550816
$ws.Range("I34").Value = 7341.278
$ws.Range("J34").Value = 333333340
$ws.Range("K34").Value = 7341.278
$ws.Range("L34").Value = 333333340
$ws.Range("M34").Value = -7139.278
$ws.Range("N34").Value = -333333744

# Row 35
$ws.Range("H35").Value = 709.0909
$ws.Range("I35").Value = 709.0909
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 709.0909
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -415.0909
$ws.Range("N35").ClearContents()

# Row 38
$ws.Range("H38").Value = 36694.668
$ws.Range("I38").Value = 30000
$ws.Range("K38").Value = 30000
$ws.Range("M38").Value = -29623

# Row 46
$ws.Range("H46").Value = 36694.668
$ws.Range("I46").Value = 30000
$ws.Range("K46").Value = 30000
$ws.Range("M46").Value = -29789

# Row 140
$ws.Range("H140").Value = 31638
$ws.Range("J140").Value = 31638
$ws.Range("L140").Value = 31638
$ws.Range("N140").Value = -41998

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 275.25
$ws.Range("I13").Value = 350.5
$ws.Range("J13").Value = 200
$ws.Range("K13").Value = 1051.5
$ws.Range("L13").Value = 600
$ws.Range("M13").Value = -883.5
$ws.Range("N13").Value = -936

# Row 99
$ws.Range("H99").Value = 2059.0588
$ws.Range("I99").Value = 1446.6666
$ws.Range("J99").Value = 2393.0908
$ws.Range("K99").Value = 4339.9998
$ws.Range("L99").Value = 7179.2724
$ws.Range("M99").Value = -2093.9998
$ws.Range("N99").Value = -11671.2724

# Row 113
$ws.Range("H113").Value = 1201.5769
$ws.Range("I113").Value = 454.91666
$ws.Range("J113").Value = 1841.5714
$ws.Range("K113").Value = 1364.74998
$ws.Range("L113").Value = 5524.7142
$ws.Range("M113").Value = 805.2500199999999
$ws.Range("N113").Value = -9864.7142

# Row 136
$ws.Range("H136").Value = 3768.7856
$ws.Range("I136").Value = 1865
$ws.Range("J136").Value = 4086.0833
$ws.Range("K136").Value = 5595
$ws.Range("L136").Value = 12258.2499
$ws.Range("M136").Value = -495
$ws.Range("N136").Value = -22458.2499

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 16669522
$ws.Range("I80").Value = 25643796
$ws.Range("J80").Value = 3014.1428
$ws.Range("K80").Value = 25643796
$ws.Range("L80").Value = 3014.1428
$ws.Range("M80").Value = -25642798
$ws.Range("N80").Value = -5010.1428

# Row 83
$ws.Range("H83").Value = 16669522
$ws.Range("I83").Value = 25643796
$ws.Range("J83").Value = 3014.1428
$ws.Range("K83").Value = 128218980
$ws.Range("L83").Value = 15070.714
$ws.Range("M83").Value = -128213988
$ws.Range("N83").Value = -25054.714

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 2609.9
$ws.Range("I82").Value = 2503.8
$ws.Range("J82").Value = 2716
$ws.Range("K82").Value = 2503.8
$ws.Range("L82").Value = 2716
$ws.Range("M82").Value = -2142.8
$ws.Range("N82").Value = -3438

# Row 85
$ws.Range("H85").Value = 2609.9
$ws.Range("I85").Value = 2503.8
$ws.Range("J85").Value = 2716
$ws.Range("K85").Value = 2503.8
$ws.Range("L85").Value = 2716
$ws.Range("M85").Value = -1255.8
$ws.Range("N85").Value = -5212

# Row 140
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360
